$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new NSE value in B4
$ws.Range("B4").Value = 0.95688099999999998

# Widen column B to fit the new numeric data (stored width ends up as 15)
$ws.Columns("B").ColumnWidth = 14.14

# Update the active selection to C4
$ws.Range("C4").Select()
